$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 19
$ws.Range("H19").Value = 556.8333
$ws.Range("J19").Value = 484.85715
$ws.Range("L19").Value = 484.85715
$ws.Range("N19").Value = -834.85715
# Row 33
$ws.Range("H33").Value = 261
$ws.Range("I33").Value = 172.5
$ws.Range("J33").Value = 1500
$ws.Range("K33").Value = 172.5
$ws.Range("L33").Value = 1500
$ws.Range("M33").Value = 56.5
$ws.Range("N33").Value = -1958
# Row 62
$ws.Range("H62").Value = 6463.4165
$ws.Range("I62").Value = 5695.875
$ws.Range("K62").Value = 5695.875
$ws.Range("M62").Value = -5071.875
# Row 65
$ws.Range("H65").Value = 6463.4165
$ws.Range("I65").Value = 5695.875
$ws.Range("K65").Value = 28479.375
$ws.Range("M65").Value = -25359.375
# Row 111
$ws.Range("H111").Value = 607.2
$ws.Range("I111").Value = 565.4286
$ws.Range("J111").Value = 704.6667
$ws.Range("K111").Value = 1696.2858
$ws.Range("L111").Value = 2114.0001
$ws.Range("M111").Value = 1370.7142
$ws.Range("N111").Value = -8248.000100000001
# Row 113
$ws.Range("H113").Value = 6321.706
$ws.Range("I113").Value = 5207.8887
$ws.Range("K113").Value = 5207.8887
$ws.Range("M113").Value = -1953.8887
# Row 115
$ws.Range("H115").Value = 65.666664
$ws.Range("I115").Value = 65.666664
$ws.Range("K115").Value = 196.999992
$ws.Range("M115").Value = 1370.000008
# Row 116
$ws.Range("H116").Value = 3821.75
$ws.Range("I116").Value = 2357.8
$ws.Range("J116").Value = 4867.4287
$ws.Range("K116").Value = 2357.8
$ws.Range("L116").Value = 4867.4287
$ws.Range("M116").Value = 1084.2
$ws.Range("N116").Value = -11751.4287
# Row 125
$ws.Range("H125").Value = 129936.375
$ws.Range("I125").Value = 1447.5
$ws.Range("J125").Value = 172766
$ws.Range("K125").Value = 13027.5
$ws.Range("L125").Value = 1554894
$ws.Range("M125").Value = -10567.5
$ws.Range("N125").Value = -1559814
# Row 132
$ws.Range("H132").Value = 3627.6667
$ws.Range("I132").Value = 5083.8335
$ws.Range("J132").Value = 715.3333
$ws.Range("K132").Value = 15251.5005
$ws.Range("L132").Value = 2145.9999
$ws.Range("M132").Value = -12721.5005
$ws.Range("N132").Value = -7205.9999
# Row 137
$ws.Range("H137").Value = 1576.7
$ws.Range("I137").Value = 1223.625
$ws.Range("K137").Value = 3670.875
$ws.Range("M137").Value = -1120.875
# Row 138
$ws.Range("H138").Value = 1903.1765
$ws.Range("J138").Value = 2722
$ws.Range("L138").Value = 8166
$ws.Range("N138").Value = -18446

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 74
$ws.Range("H74").Value = 4331.25
$ws.Range("I74").Value = 4337
$ws.Range("J74").Value = 4314
$ws.Range("K74").Value = 4337
$ws.Range("L74").Value = 4314
$ws.Range("M74").Value = -3463
$ws.Range("N74").Value = -6062
# Row 77
$ws.Range("H77").Value = 4331.25
$ws.Range("I77").Value = 4337
$ws.Range("J77").Value = 4314
$ws.Range("K77").Value = 21685
$ws.Range("L77").Value = 21570
$ws.Range("M77").Value = -17317
$ws.Range("N77").Value = -30306
# Row 97
$ws.Range("H97").Value = 644.625
$ws.Range("I97").Value = 593.8570999999999
$ws.Range("K97").Value = 593.8570999999999
$ws.Range("M97").Value = -97.85709999999995
# Row 110
$ws.Range("H110").Value = 2672.7778
$ws.Range("I110").Value = 1411.2
$ws.Range("K110").Value = 1411.2
$ws.Range("M110").Value = 633.8
# Row 132
$ws.Range("H132").Value = 4751.625
$ws.Range("I132").Value = 4787.5713
$ws.Range("K132").Value = 14362.7139
$ws.Range("M132").Value = -11832.7139

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 41
$ws.Range("H41").Value = 5686.3335
$ws.Range("J41").Value = 10000
$ws.Range("L41").Value = 10000
$ws.Range("N41").Value = -10856
# Row 47
$ws.Range("H47").Value = 123075
$ws.Range("I47").Value = 12000
$ws.Range("K47").Value = 12000
$ws.Range("M47").Value = -11434
# Row 58
$ws.Range("H58").Value = 2364.1667
$ws.Range("I58").Value = 1460.7
$ws.Range("K58").Value = 1460.7
$ws.Range("M58").Value = -1257.7
# Row 59
$ws.Range("H59").Value = 27350.3
$ws.Range("I59").Value = 17834.666
$ws.Range("J59").Value = 31428.428
$ws.Range("K59").Value = 17834.666
$ws.Range("L59").Value = 31428.428
$ws.Range("M59").Value = -16689.666
$ws.Range("N59").Value = -33718.428
# Row 60
$ws.Range("H60").Value = 21189.1
$ws.Range("I60").Value = 21981.834
$ws.Range("K60").Value = 21981.834
$ws.Range("M60").Value = -21470.834
# Row 136
$ws.Range("H136").Value = 2364.1667
$ws.Range("I136").Value = 1460.7
$ws.Range("K136").Value = 4382.1
$ws.Range("M136").Value = -1832.1

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 137
$ws.Range("H137").Value = 1999.6666
$ws.Range("J137").Value = 1999.5
$ws.Range("L137").Value = 5998.5
$ws.Range("N137").Value = -16198.5

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 2
$ws.Range("H2").Value = 170.6923
$ws.Range("I2").Value = 252.42857
$ws.Range("J2").Value = 75.333336
$ws.Range("K2").Value = 252.42857
$ws.Range("L2").Value = 75.333336
$ws.Range("M2").Value = -139.42857
$ws.Range("N2").Value = -301.333336
# Row 70
$ws.Range("H70").Value = 2339.4
$ws.Range("I70").Value = 2399
$ws.Range("J70").Value = 2250
$ws.Range("K70").Value = 2399
$ws.Range("L70").Value = 2250
$ws.Range("M70").Value = -2129
$ws.Range("N70").Value = -2790
# Row 73
$ws.Range("H73").Value = 2339.4
$ws.Range("I73").Value = 2399
$ws.Range("J73").Value = 2250
$ws.Range("K73").Value = 2399
$ws.Range("L73").Value = 2250
$ws.Range("M73").Value = -1463
$ws.Range("N73").Value = -4122
# Row 80
$ws.Range("H80").Value = 2379.5
$ws.Range("I80").Value = 2279.8
$ws.Range("J80").Value = 2479.2
$ws.Range("K80").Value = 2279.8
$ws.Range("L80").Value = 2479.2
$ws.Range("M80").Value = -1281.8
$ws.Range("N80").Value = -4475.2
# Row 83
$ws.Range("H83").Value = 2379.5
$ws.Range("I83").Value = 2279.8
$ws.Range("J83").Value = 2479.2
$ws.Range("K83").Value = 11399
$ws.Range("L83").Value = 12396
$ws.Range("M83").Value = -6407
$ws.Range("N83").Value = -22380
# Row 97
$ws.Range("H97").Value = 570.2222
$ws.Range("I97").Value = 581.625
$ws.Range("J97").Value = 479
$ws.Range("K97").Value = 581.625
$ws.Range("L97").Value = 479
$ws.Range("M97").Value = -85.625
$ws.Range("N97").Value = -1471
# Row 114
$ws.Range("H114").Value = 312722
$ws.Range("J114").Value = 312722
$ws.Range("L114").Value = 312722
$ws.Range("N114").Value = -321400
# Row 126
$ws.Range("H126").Value = 12417.5
$ws.Range("I126").Value = 9335.333000000001
$ws.Range("J126").Value = 15499.667
$ws.Range("K126").Value = 28005.999
$ws.Range("L126").Value = 46499.001
$ws.Range("M126").Value = -25535.999
$ws.Range("N126").Value = -51439.001

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 1200
$ws.Range("I22").Value = 1000
$ws.Range("J22").Value = 2000
$ws.Range("K22").Value = 1000
$ws.Range("L22").Value = 2000
$ws.Range("M22").Value = -705
$ws.Range("N22").Value = -2590
# Row 27
$ws.Range("H27").Value = 1200
$ws.Range("I27").Value = 1000
$ws.Range("J27").Value = 2000
$ws.Range("K27").Value = 1000
$ws.Range("L27").Value = 2000
$ws.Range("M27").Value = -893
$ws.Range("N27").Value = -2214
# Row 45
$ws.Range("H45").Value = 0
$ws.Range("I45").Value = 0
$ws.Range("K45").Value = 0
$ws.Range("M45").ClearContents()
# Row 48
$ws.Range("H48").Value = 20000
$ws.Range("J48").Value = 20000
$ws.Range("L48").Value = 20000
$ws.Range("N48").Value = -21322
# Row 61
$ws.Range("H61").Value = 1697
$ws.Range("I61").Value = 1742.1818
$ws.Range("K61").Value = 1742.1818
$ws.Range("M61").Value = -1540.1818
# Row 113
$ws.Range("H113").Value = 1697
$ws.Range("I113").Value = 1742.1818
$ws.Range("K113").Value = 1742.1818
$ws.Range("M113").Value = 427.8181999999999

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 18
$ws.Range("H18").Value = 1000
$ws.Range("I18").Value = 1000
$ws.Range("K18").Value = 1000
$ws.Range("M18").Value = -827
# Row 81
$ws.Range("H81").Value = 552.8570999999999
$ws.Range("I81").Value = 561.6667
$ws.Range("K81").Value = 1123.3334
$ws.Range("M81").Value = -62.33339999999998
# Row 84
$ws.Range("H84").Value = 552.8570999999999
$ws.Range("I84").Value = 561.6667
$ws.Range("K84").Value = 5616.666999999999
$ws.Range("M84").Value = -312.6669999999995
# Row 96
$ws.Range("H96").Value = 969.75
$ws.Range("I96").Value = 969.75
$ws.Range("K96").Value = 969.75
$ws.Range("M96").Value = 403.25
# Row 107
$ws.Range("H107").Value = 1906.2222
$ws.Range("I107").Value = 2101.7
$ws.Range("J107").Value = 1661.875
$ws.Range("K107").Value = 6305.099999999999
$ws.Range("L107").Value = 4985.625
$ws.Range("M107").Value = -4385.099999999999
$ws.Range("N107").Value = -8825.625
